$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5608
$ws.Range("K3").Value = 5755
$ws.Range("E4").Value = 2033
$ws.Range("F4").Value = 1913
$ws.Range("H4").Value = 1740
$ws.Range("K4").Value = 1195
$ws.Range("K5").Value = 411
$ws.Range("K6").Value = 6389
$ws.Range("E7").Value = 26039
$ws.Range("F7").Value = 24106
$ws.Range("H7").Value = 26053
$ws.Range("K7").Value = 19358

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 354
$ws.Range("K3").Value = 389
$ws.Range("K6").Value = 430
$ws.Range("K7").Value = 1282

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 155
$ws.Range("K7").Value = 431

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 303
$ws.Range("K6").Value = 244
$ws.Range("K7").Value = 830

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 114
$ws.Range("K7").Value = 327

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 216
$ws.Range("K5").Value = 29
$ws.Range("K7").Value = 656

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 111
$ws.Range("K4").Value = 17
$ws.Range("K7").Value = 445

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 133
$ws.Range("K7").Value = 325

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 172
$ws.Range("K6").Value = 140
$ws.Range("K7").Value = 567
$ws.Range("K8").Value = 1282
$ws.Range("K11").Value = 367
$ws.Range("K14").Value = 101
$ws.Range("K15").Value = 198
$ws.Range("K19").Value = 567
$ws.Range("K20").Value = 454
$ws.Range("K23").Value = 199
$ws.Range("K25").Value = 89
$ws.Range("K29").Value = 1049
$ws.Range("K31").Value = 209
$ws.Range("K33").Value = 830
$ws.Range("K36").Value = 252
$ws.Range("K37").Value = 656
$ws.Range("K42").Value = 719
$ws.Range("K45").Value = 23
$ws.Range("K47").Value = 135
$ws.Range("K50").Value = 93
$ws.Range("K51").Value = 247
$ws.Range("K52").Value = 502
$ws.Range("K54").Value = 371
$ws.Range("K55").Value = 213
$ws.Range("E63").Value = 369
$ws.Range("F63").Value = 197
$ws.Range("H63").Value = 291
$ws.Range("K63").Value = 55
$ws.Range("K64").Value = 125
$ws.Range("K65").Value = 445
$ws.Range("K67").Value = 744
$ws.Range("K68").Value = 50
$ws.Range("K72").Value = 90
$ws.Range("K73").Value = 171
$ws.Range("K76").Value = 264
$ws.Range("K77").Value = 132
$ws.Range("K79").Value = 488
$ws.Range("K83").Value = 431
$ws.Range("K85").Value = 911
$ws.Range("K87").Value = 32
$ws.Range("K89").Value = 283
$ws.Range("K90").Value = 176
$ws.Range("K91").Value = 216
$ws.Range("K94").Value = 261
$ws.Range("K95").Value = 327
$ws.Range("K96").Value = 207
$ws.Range("K97").Value = 155
$ws.Range("K99").Value = 325
$ws.Range("E101").Value = 26039
$ws.Range("F101").Value = 24106
$ws.Range("H101").Value = 26053
$ws.Range("K101").Value = 19358

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 209

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 214
$ws.Range("K3").Value = 267
$ws.Range("K6").Value = 206
$ws.Range("K7").Value = 744

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 198
$ws.Range("K7").Value = 371

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 299
$ws.Range("K7").Value = 1049

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 176
$ws.Range("K7").Value = 567

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 138
$ws.Range("K7").Value = 264

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 53
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 220
$ws.Range("K7").Value = 719

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 65
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 213

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 56
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 207

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 102
$ws.Range("K7").Value = 216

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 163
$ws.Range("K7").Value = 488

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 148
$ws.Range("K3").Value = 148
$ws.Range("K6").Value = 130
$ws.Range("K7").Value = 454

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 73
$ws.Range("K7").Value = 252

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 183
$ws.Range("K7").Value = 567

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 261

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 38
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 70
$ws.Range("K7").Value = 198

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K2").Value = 24
$ws.Range("K3").Value = 13
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 126
$ws.Range("K7").Value = 367

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 58
$ws.Range("K7").Value = 171

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 172

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K3").Value = 31
$ws.Range("K7").Value = 155

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 79
$ws.Range("K7").Value = 283

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 247

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K2").Value = 22
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 297
$ws.Range("K3").Value = 308
$ws.Range("K7").Value = 911

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K2").Value = 4
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 136
$ws.Range("K7").Value = 502

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 32
